$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated result values (columns B,C,D,E,F,H,J,L,M) for rows 2-25
# per commit 'case with 380 kV done'

# Row 2 (A2 = 0)
$ws.Cells.Item(2, 2).Value = 3.276109556179165
$ws.Cells.Item(2, 3).Value = 0.1710648090741813
$ws.Cells.Item(2, 4).Value = 0.0508899230681692
$ws.Cells.Item(2, 5).Value = 0.04304938155374316
$ws.Cells.Item(2, 6).Value = 6.946608703946481
$ws.Cells.Item(2, 8).Value = 0.07973214163530429
$ws.Cells.Item(2, 10).Value = 0.1678877235440694
$ws.Cells.Item(2, 12).Value = 0.2796376590313159
$ws.Cells.Item(2, 13).Value = 0.6080991648005281

# Row 3 (A3 = 1)
$ws.Cells.Item(3, 2).Value = 3.240473125576557
$ws.Cells.Item(3, 3).Value = 0.1593196973515489
$ws.Cells.Item(3, 4).Value = 0.04484173880611309
$ws.Cells.Item(3, 5).Value = 0.04265931115019939
$ws.Cells.Item(3, 6).Value = 6.77616795570583
$ws.Cells.Item(3, 8).Value = 0.07973214163530429
$ws.Cells.Item(3, 10).Value = 0.1663173386366132
$ws.Cells.Item(3, 12).Value = 0.2815785064646192
$ws.Cells.Item(3, 13).Value = 0.6051180434261525

# Row 4 (A4 = 2)
$ws.Cells.Item(4, 2).Value = 3.221219722281205
$ws.Cells.Item(4, 3).Value = 0.1523128592112073
$ws.Cells.Item(4, 4).Value = 0.04112014866902314
$ws.Cells.Item(4, 5).Value = 0.0424142896310622
$ws.Cells.Item(4, 6).Value = 6.67300045888274
$ws.Cells.Item(4, 8).Value = 0.07973214163530429
$ws.Cells.Item(4, 10).Value = 0.1653364572585083
$ws.Cells.Item(4, 12).Value = 0.2829349338570992
$ws.Cells.Item(4, 13).Value = 0.6037491444177547

# Row 5 (A5 = 3)
$ws.Cells.Item(5, 2).Value = 3.214033277698775
$ws.Cells.Item(5, 3).Value = 0.1495084758911673
$ws.Cells.Item(5, 4).Value = 0.03960130796437511
$ws.Cells.Item(5, 5).Value = 0.04231303790321927
$ws.Cells.Item(5, 6).Value = 6.631325723878859
$ws.Cells.Item(5, 8).Value = 0.07973214163530429
$ws.Cells.Item(5, 10).Value = 0.164932456785035
$ws.Cells.Item(5, 12).Value = 0.2835291633455128
$ws.Cells.Item(5, 13).Value = 0.6033072436213018

# Row 6 (A6 = 4)
$ws.Cells.Item(6, 2).Value = 3.212879769716096
$ws.Cells.Item(6, 3).Value = 0.1490458731746003
$ws.Cells.Item(6, 4).Value = 0.03934896152247802
$ws.Cells.Item(6, 5).Value = 0.04229613989224124
$ws.Cells.Item(6, 6).Value = 6.624427641421022
$ws.Cells.Item(6, 8).Value = 0.07973214163530429
$ws.Cells.Item(6, 10).Value = 0.164865111358047
$ws.Cells.Item(6, 12).Value = 0.2836303411314489
$ws.Cells.Item(6, 13).Value = 0.6032408653933388

# Row 7 (A7 = 5)
$ws.Cells.Item(7, 2).Value = 3.221120134760326
$ws.Cells.Item(7, 3).Value = 0.1522748326541432
$ws.Cells.Item(7, 4).Value = 0.04109967455033825
$ws.Cells.Item(7, 5).Value = 0.04241292981897615
$ws.Cells.Item(7, 6).Value = 6.672436941662795
$ws.Cells.Item(7, 8).Value = 0.07973214163530429
$ws.Cells.Item(7, 10).Value = 0.1653310262346643
$ws.Cells.Item(7, 12).Value = 0.2829427798471542
$ws.Cells.Item(7, 13).Value = 0.6037427155094051

# Row 8 (A8 = 6)
$ws.Cells.Item(8, 2).Value = 3.26327599087881
$ws.Cells.Item(8, 3).Value = 0.1669722901542912
$ws.Cells.Item(8, 4).Value = 0.04880600971324611
$ws.Cells.Item(8, 5).Value = 0.042916020388021
$ws.Cells.Item(8, 6).Value = 6.887528931128372
$ws.Cells.Item(8, 8).Value = 0.07973214163530429
$ws.Cells.Item(8, 10).Value = 0.1673496551617113
$ws.Cells.Item(8, 12).Value = 0.280272700262195
$ws.Cells.Item(8, 13).Value = 0.6069754091130477

# Row 9 (A9 = 7)
$ws.Cells.Item(9, 2).Value = 3.366860676287672
$ws.Cells.Item(9, 3).Value = 0.1974436504785331
$ws.Cells.Item(9, 4).Value = 0.06386695805304043
$ws.Cells.Item(9, 5).Value = 0.04385957043913535
$ws.Cells.Item(9, 6).Value = 7.321403551961737
$ws.Cells.Item(9, 8).Value = 0.07973214163530429
$ws.Cells.Item(9, 10).Value = 0.1711803682580424
$ws.Cells.Item(9, 12).Value = 0.2763419084387948
$ws.Cells.Item(9, 13).Value = 0.6169838819495723

# Row 10 (A10 = 8)
$ws.Cells.Item(10, 2).Value = 3.455831124156191
$ws.Cells.Item(10, 3).Value = 0.2208750269706172
$ws.Cells.Item(10, 4).Value = 0.07491882383139625
$ws.Cells.Item(10, 5).Value = 0.04452767904979638
$ws.Cells.Item(10, 6).Value = 7.647994882885371
$ws.Cells.Item(10, 8).Value = 0.07973214163530429
$ws.Cells.Item(10, 10).Value = 0.1739230815191277
$ws.Cells.Item(10, 12).Value = 0.2742474856621655
$ws.Cells.Item(10, 13).Value = 0.626586735843496

# Row 11 (A11 = 9)
$ws.Cells.Item(11, 2).Value = 3.499127166959227
$ws.Cells.Item(11, 3).Value = 0.2317699354061631
$ws.Cells.Item(11, 4).Value = 0.07994760146434032
$ws.Cells.Item(11, 5).Value = 0.04482640719533126
$ws.Cells.Item(11, 6).Value = 7.798370922274842
$ws.Cells.Item(11, 8).Value = 0.07973214163530429
$ws.Cells.Item(11, 10).Value = 0.1751565949932612
$ws.Cells.Item(11, 12).Value = 0.2734665935958844
$ws.Cells.Item(11, 13).Value = 0.6314469276339381

# Row 12 (A12 = 10)
$ws.Cells.Item(12, 2).Value = 3.515930204245251
$ws.Cells.Item(12, 3).Value = 0.2359301391719839
$ws.Cells.Item(12, 4).Value = 0.08185234598535374
$ws.Cells.Item(12, 5).Value = 0.04493879929563516
$ws.Cells.Item(12, 6).Value = 7.855582152492218
$ws.Cells.Item(12, 8).Value = 0.07973214163530429
$ws.Cells.Item(12, 10).Value = 0.1756217652978904
$ws.Cells.Item(12, 12).Value = 0.2731955732921065
$ws.Cells.Item(12, 13).Value = 0.6333583062042791

# Row 13 (A13 = 11)
$ws.Cells.Item(13, 2).Value = 3.51229320403462
$ws.Cells.Item(13, 3).Value = 0.2350326186821121
$ws.Cells.Item(13, 4).Value = 0.08144210043977296
$ws.Cells.Item(13, 5).Value = 0.04491462581875361
$ws.Cells.Item(13, 6).Value = 7.843248697420904
$ws.Cells.Item(13, 8).Value = 0.07973214163530429
$ws.Cells.Item(13, 10).Value = 0.1755216669565165
$ws.Cells.Item(13, 12).Value = 0.2732528448448761
$ws.Cells.Item(13, 13).Value = 0.6329434985760045

# Row 14 (A14 = 12)
$ws.Cells.Item(14, 2).Value = 3.50050137995737
$ws.Cells.Item(14, 3).Value = 0.2321115020560001
$ws.Cells.Item(14, 4).Value = 0.08010429520381024
$ws.Cells.Item(14, 5).Value = 0.04483566828463115
$ws.Cells.Item(14, 6).Value = 7.803072330737677
$ws.Cells.Item(14, 8).Value = 0.07973214163530429
$ws.Cells.Item(14, 10).Value = 0.1751949030176227
$ws.Cells.Item(14, 12).Value = 0.2734438019901333
$ws.Cells.Item(14, 13).Value = 0.6316027551318584

# Row 15 (A15 = 13)
$ws.Cells.Item(15, 2).Value = 3.493331706005335
$ws.Cells.Item(15, 3).Value = 0.230326750419664
$ws.Cells.Item(15, 4).Value = 0.07928491852527486
$ws.Cells.Item(15, 5).Value = 0.04478721000650054
$ws.Cells.Item(15, 6).Value = 7.778498141296836
$ws.Cells.Item(15, 8).Value = 0.07973214163530429
$ws.Cells.Item(15, 10).Value = 0.1749945017808585
$ws.Cells.Item(15, 12).Value = 0.2735639828780592
$ws.Cells.Item(15, 13).Value = 0.6307907545528622

# Row 16 (A16 = 14)
$ws.Cells.Item(16, 2).Value = 3.453058585346923
$ws.Cells.Item(16, 3).Value = 0.2201678180024373
$ws.Cells.Item(16, 4).Value = 0.07459022988484776
$ws.Cells.Item(16, 5).Value = 0.0445080533753428
$ws.Cells.Item(16, 6).Value = 7.638204441119342
$ws.Cells.Item(16, 8).Value = 0.07973214163530429
$ws.Cells.Item(16, 10).Value = 0.1738421922476121
$ws.Cells.Item(16, 12).Value = 0.274301974807706
$ws.Cells.Item(16, 13).Value = 0.6262790248025567

# Row 17 (A17 = 15)
$ws.Cells.Item(17, 2).Value = 3.429076630071847
$ws.Cells.Item(17, 3).Value = 0.2139964889551038
$ws.Cells.Item(17, 4).Value = 0.07171070266721813
$ws.Cells.Item(17, 5).Value = 0.04433548104457419
$ws.Cells.Item(17, 6).Value = 7.552606890770249
$ws.Cells.Item(17, 8).Value = 0.07973214163530429
$ws.Cells.Item(17, 10).Value = 0.1731317329668336
$ws.Cells.Item(17, 12).Value = 0.2747987083805299
$ws.Cells.Item(17, 13).Value = 0.6236373323799569

# Row 18 (A18 = 16)
$ws.Cells.Item(18, 2).Value = 3.415548437924372
$ws.Cells.Item(18, 3).Value = 0.2104690656431956
$ws.Cells.Item(18, 4).Value = 0.07005457355246847
$ws.Cells.Item(18, 5).Value = 0.04423573169094563
$ws.Cells.Item(18, 6).Value = 7.503543190871682
$ws.Cells.Item(18, 8).Value = 0.07973214163530429
$ws.Cells.Item(18, 10).Value = 0.1727217565811934
$ws.Cells.Item(18, 12).Value = 0.2751005952800014
$ws.Cells.Item(18, 13).Value = 0.6221641786322962

# Row 19 (A19 = 17)
$ws.Cells.Item(19, 2).Value = 3.411013586060733
$ws.Cells.Item(19, 3).Value = 0.2092785283327032
$ws.Cells.Item(19, 4).Value = 0.06949384520501667
$ws.Cells.Item(19, 5).Value = 0.04420187352746918
$ws.Cells.Item(19, 6).Value = 7.486960004764939
$ws.Cells.Item(19, 8).Value = 0.07973214163530429
$ws.Cells.Item(19, 10).Value = 0.1725827125614465
$ws.Cells.Item(19, 12).Value = 0.2752055886384639
$ws.Cells.Item(19, 13).Value = 0.6216733361923659

# Row 20 (A20 = 18)
$ws.Cells.Item(20, 2).Value = 3.43160204915938
$ws.Cells.Item(20, 3).Value = 0.2146511393309822
$ws.Cells.Item(20, 4).Value = 0.0720172210202179
$ws.Cells.Item(20, 5).Value = 0.04435390228202252
$ws.Cells.Item(20, 6).Value = 7.561701276874913
$ws.Cells.Item(20, 8).Value = 0.07973214163530429
$ws.Cells.Item(20, 10).Value = 0.1732075004164493
$ws.Cells.Item(20, 12).Value = 0.274744156004175
$ws.Cells.Item(20, 13).Value = 0.6239137539948416

# Row 21 (A21 = 19)
$ws.Cells.Item(21, 2).Value = 3.503953842904082
$ws.Cells.Item(21, 3).Value = 0.2329685624219167
$ws.Cells.Item(21, 4).Value = 0.08049722642998347
$ws.Cells.Item(21, 5).Value = 0.04485887968984148
$ws.Cells.Item(21, 6).Value = 7.814865797725247
$ws.Cells.Item(21, 8).Value = 0.07973214163530429
$ws.Cells.Item(21, 10).Value = 0.1752909331928869
$ws.Cells.Item(21, 12).Value = 0.2733870434741732
$ws.Cells.Item(21, 13).Value = 0.6319946372562484

# Row 22 (A22 = 20)
$ws.Cells.Item(22, 2).Value = 3.553617598480571
$ws.Cells.Item(22, 3).Value = 0.2451416730317533
$ws.Cells.Item(22, 4).Value = 0.08604218711603551
$ws.Cells.Item(22, 5).Value = 0.04518466635447282
$ws.Cells.Item(22, 6).Value = 7.981883610621992
$ws.Cells.Item(22, 8).Value = 0.07973214163530429
$ws.Cells.Item(22, 10).Value = 0.1766413492328702
$ws.Cells.Item(22, 12).Value = 0.2726439728954659
$ws.Cells.Item(22, 13).Value = 0.6376894978826968

# Row 23 (A23 = 21)
$ws.Cells.Item(23, 2).Value = 3.526892934420403
$ws.Cells.Item(23, 3).Value = 0.2386259974128677
$ws.Cells.Item(23, 4).Value = 0.08308238932497147
$ws.Cells.Item(23, 5).Value = 0.04501117042341196
$ws.Cells.Item(23, 6).Value = 7.892597931706518
$ws.Cells.Item(23, 8).Value = 0.07973214163530429
$ws.Cells.Item(23, 10).Value = 0.1759216002456903
$ws.Cells.Item(23, 12).Value = 0.2730274072280281
$ws.Cells.Item(23, 13).Value = 0.6346121329801733

# Row 24 (A24 = 22)
$ws.Cells.Item(24, 2).Value = 3.430459499499534
$ws.Cells.Item(24, 3).Value = 0.2143551079150825
$ws.Cells.Item(24, 4).Value = 0.07187864612041039
$ws.Cells.Item(24, 5).Value = 0.04434557570901987
$ws.Cells.Item(24, 6).Value = 7.557589246109785
$ws.Cells.Item(24, 8).Value = 0.07973214163530429
$ws.Cells.Item(24, 10).Value = 0.1731732507063981
$ws.Cells.Item(24, 12).Value = 0.2747687683350222
$ws.Cells.Item(24, 13).Value = 0.623788641794313

# Row 25 (A25 = 23)
$ws.Cells.Item(25, 2).Value = 3.336587318046099
$ws.Cells.Item(25, 3).Value = 0.1890198127254052
$ws.Cells.Item(25, 4).Value = 0.05979622062119461
$ws.Cells.Item(25, 5).Value = 0.04360880639189979
$ws.Cells.Item(25, 6).Value = 7.202688595208883
$ws.Cells.Item(25, 8).Value = 0.07973214163530429
$ws.Cells.Item(25, 10).Value = 0.1701571032338549
$ws.Cells.Item(25, 12).Value = 0.2772657915718284
$ws.Cells.Item(25, 13).Value = 0.6138822324665298
